$d = $word.ActiveDocument

# The page footer (injected by the static-site generator that produced this
# document) consists of three consecutive paragraphs right before the
# trailing blank paragraph / page-break paragraph at the end of the body:
#   1. an empty paragraph
#   2. "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3. "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#       pages. Original theme under Creative Commons Attribution"
# The site was rebuilt and this footer block is no longer present, so all
# three paragraphs must be removed while leaving everything else intact.

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $target = $i
        break
    }
}

if ($null -ne $target) {
    # Include the blank paragraph immediately before the "Ver no Jupiter..."
    # paragraph and the copyright paragraph immediately after it.
    $startPara = $d.Paragraphs.Item($target - 1)
    $endPara   = $d.Paragraphs.Item($target + 1)
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
